# Trade #3 closed at 2026-02-17 04:05:53 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99
$summary.Range("B4").Value = -0.01
$summary.Range("B5").Value = -0.07000000000000001
$summary.Range("B6").Value = 3
$summary.Range("B7").Value = 1
$summary.Range("B9").Value = 33.33

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999
$status.Range("D4").Value = 3
$status.Range("E4").Value = -0.01
$status.Range("F4").Value = -0.01
$status.Range("G4").Value = 33.33

# ---- New trade row data (Trade #3) ----
$tradeNum = 3
$tradeDate = "2026-02-17"
$tradeTime = "04:05:47"
$tradeStrategy = "MarketMaking"
$tradeSide = "DOWN"
$entryPrice = 0.73
$exitPrice = 0.76
$tradeStatus = "CLOSED"
$pnlPct = 4.1096
$pnlDollar = 0.03
$capitalAfter = 99.98999999999999
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$durationMin = 0.12

# ---- All Trades sheet: append row 4 ----
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A4").Value = $tradeNum
# Force the date-looking string to stay plain text (not auto-converted to a
# date serial) the way the original file stores it, then strip the
# temporary Text number format so no stray style sticks to the cell.
$allTrades.Range("B4").NumberFormat = "@"
$allTrades.Range("B4").Value = $tradeDate
$allTrades.Range("B4").ClearFormats()
$allTrades.Range("C4").Value = $tradeTime
$allTrades.Range("D4").Value = $tradeStrategy
$allTrades.Range("E4").Value = $tradeSide
$allTrades.Range("F4").Value = $entryPrice
$allTrades.Range("G4").Value = $exitPrice
$allTrades.Range("H4").Value = $tradeStatus
$allTrades.Range("I4").Value = $pnlPct
$allTrades.Range("J4").Value = $pnlDollar
$allTrades.Range("K4").Value = $capitalAfter
$allTrades.Range("L4").Value = $entrySlippage
$allTrades.Range("M4").Value = $exitSlippage
$allTrades.Range("N4").Value = $confidence
$allTrades.Range("O4").Value = $entryReason
$allTrades.Range("P4").Value = $exitReason
$allTrades.Range("Q4").Value = $durationMin

# ---- MarketMaking sheet: append row 4 ----
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A4").Value = $tradeNum
$mm.Range("B4").NumberFormat = "@"
$mm.Range("B4").Value = $tradeDate
$mm.Range("B4").ClearFormats()
$mm.Range("C4").Value = $tradeTime
$mm.Range("D4").Value = $tradeStrategy
$mm.Range("E4").Value = $tradeSide
$mm.Range("F4").Value = $entryPrice
$mm.Range("G4").Value = $exitPrice
$mm.Range("H4").Value = $tradeStatus
$mm.Range("I4").Value = $pnlPct
$mm.Range("J4").Value = $pnlDollar
$mm.Range("K4").Value = $capitalAfter
$mm.Range("L4").Value = $entrySlippage
$mm.Range("M4").Value = $exitSlippage
$mm.Range("N4").Value = $confidence
$mm.Range("O4").Value = $entryReason
$mm.Range("P4").Value = $exitReason
$mm.Range("Q4").Value = $durationMin
